$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.6
$ws.Range("I2").Value = 5.75
$ws.Range("AH2").Value = 15
$ws.Range("AI2").Value = 29
$ws.Range("AS2").Value = 126

# Row 4
$ws.Range("I4").Value = 6
$ws.Range("K4").Value = 1.83
$ws.Range("L4").Value = 6.5
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5
$ws.Range("O4").Value = 1.62
$ws.Range("P4").Value = 2.2
$ws.Range("Q4").Value = 3.1
$ws.Range("R4").Value = 1.36
$ws.Range("S4").Value = 1.67
$ws.Range("T4").Value = 2.1
$ws.Range("U4").Value = 2.63
$ws.Range("V4").Value = 1.44
$ws.Range("W4").Value = 4.5
$ws.Range("AC4").Value = 5
$ws.Range("AE4").Value = 26
$ws.Range("AF4").Value = 126
$ws.Range("AM4").Value = 81
$ws.Range("AT4").Value = 2.1
$ws.Range("AW4").Value = 7
$ws.Range("AY4").Value = 51
$ws.Range("BA4").Value = 251

# Row 5
$ws.Range("Q5").Value = 2.03
$ws.Range("R5").Value = 1.83

# Row 6
$ws.Range("I6").Value = 7.5
$ws.Range("K6").Value = 2.6
$ws.Range("L6").Value = 6.5
$ws.Range("O6").Value = 1.15
$ws.Range("P6").Value = 5
$ws.Range("Q6").Value = 1.47
$ws.Range("R6").Value = 2.32
$ws.Range("T6").Value = 3.75
$ws.Range("U6").Value = 1.74
$ws.Range("V6").Value = 2.02
$ws.Range("W6").Value = 7.3
$ws.Range("X6").Value = 6.1
$ws.Range("Z6").Value = 7.2
$ws.Range("AC6").Value = 16
$ws.Range("AD6").Value = 9
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 19
$ws.Range("AI6").Value = 40
$ws.Range("AK6").Value = 120
$ws.Range("AL6").Value = 60
$ws.Range("AM6").Value = 50
$ws.Range("AP6").Value = 14
$ws.Range("AT6").Value = 3.55
$ws.Range("AV6").Value = 65
$ws.Range("AX6").Value = 40
$ws.Range("BA6").Value = 250

# Row 7
$ws.Range("H7").Value = 2.95
$ws.Range("I7").Value = 2.4
$ws.Range("J7").Value = 3.45
$ws.Range("K7").Value = 2.02
$ws.Range("L7").Value = 3
$ws.Range("O7").Value = 1.38
$ws.Range("P7").Value = 2.57
$ws.Range("Q7").Value = 2.12
$ws.Range("S7").Value = 1.4
$ws.Range("T7").Value = 2.52
$ws.Range("U7").Value = 1.82
$ws.Range("V7").Value = 1.78
$ws.Range("X7").Value = 15
$ws.Range("Y7").Value = 10.75
$ws.Range("AC7").Value = 7.5
$ws.Range("AD7").Value = 5.8
$ws.Range("AE7").Value = 14.5
$ws.Range("AF7").Value = 80
$ws.Range("AG7").Value = 700
$ws.Range("AH7").Value = 6.8
$ws.Range("AI7").Value = 11
$ws.Range("AJ7").Value = 9.5
$ws.Range("AK7").Value = 26
$ws.Range("AL7").Value = 22
$ws.Range("AM7").Value = 35
$ws.Range("AN7").Value = 4.85
$ws.Range("AO7").Value = 16
$ws.Range("AT7").Value = 2.47
$ws.Range("AU7").Value = 6.7
$ws.Range("AV7").Value = 60
$ws.Range("AW7").Value = 4.25
$ws.Range("AX7").Value = 12.5
$ws.Range("AY7").Value = 21

# Row 8
$ws.Range("BC8").Value = 151

# Row 11
$ws.Range("I11").Value = 4.2
$ws.Range("M11").Value = 1.1
$ws.Range("N11").Value = 6
$ws.Range("O11").Value = 1.5
$ws.Range("P11").Value = 2.37
$ws.Range("X11").Value = 7.5
$ws.Range("AE11").Value = 23
$ws.Range("AF11").Value = 101
$ws.Range("AU11").Value = 10

# Row 12
$ws.Range("M12").Value = 1.05
$ws.Range("O12").Value = 1.41
$ws.Range("P12").Value = 2.62

# Row 13
$ws.Range("G13").Value = 1.3
$ws.Range("H13").Value = 4.75
$ws.Range("I13").Value = 11
$ws.Range("L13").Value = 8.5
$ws.Range("N13").Value = 13
$ws.Range("O13").Value = 1.17
$ws.Range("U13").Value = 2.1
$ws.Range("V13").Value = 1.67
$ws.Range("X13").Value = 6.5
$ws.Range("Y13").Value = 9
$ws.Range("Z13").Value = 8
$ws.Range("AD13").Value = 9.5
$ws.Range("AE13").Value = 21
$ws.Range("AK13").Value = 126
$ws.Range("AL13").Value = 67
$ws.Range("AS13").Value = 126
$ws.Range("AZ13").Value = 201
$ws.Range("BA13").Value = 201
